# No-op edit: the underlying workbook content is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
